$d = $word.ActiveDocument

# Locate the existing date text "18/05/2020" in the document.
$rng = $d.Content
$found = $rng.Find.Execute("18/05/2020", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

$start = $rng.Start

# Replace the "18" prefix with "19", leaving "/05/2020" untouched.
$dayRange = $d.Range($start, $start + 2)
$dayRange.Text = "19"

# Force Word to keep this edited prefix as its own run (matching how Word
# really splits a run when only part of its text is retyped) while ensuring
# the run's formatting ends up identical to its neighbour, so the two runs
# carry the exact same rPr as in the target document.
$dayRange2 = $d.Range($start, $start + 2)
$dayRange2.Font.Size = 11
$dayRange3 = $d.Range($start, $start + 2)
$dayRange3.Font.Size = 12
